$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "addCategory": the two sample parentCat values now point at the
# existing "Hai Category 1" / "Hai Category 2" categories instead of the
# (now unused) "Toy" / "Hot Categories" placeholders.
# ---------------------------------------------------------------------------
$wsCat = $wb.Worksheets.Item("addCategory")
$wsCat.Range("B2").Value = "Hai Category 1"
$wsCat.Range("B3").Value = "Hai Category 2"
$null = $wsCat.Range("D9").Select()

# ---------------------------------------------------------------------------
# Sheet "addProduct": add a "unit" column (after brandName) and a
# "discountType" column (after discount, before quantity).
# ---------------------------------------------------------------------------
$wsProd = $wb.Worksheets.Item("addProduct")

# Insert the new "unit" column at D (pushes minPurchaseQty..quantity right).
$wsProd.Columns("D").Insert()
# Insert the new "discountType" column at I (old H/quantity moves to J).
$wsProd.Columns("I").Insert()

# Column widths for the newly-created columns (match their neighbours: C's
# width for D, H's width for I). The engine's ColumnWidth setter quantises
# to whole pixels, so the inputs below are chosen to land as close as
# possible to the neighbour's stored width.
$wsProd.Columns("D").ColumnWidth = 15.666666666666666
$wsProd.Columns("I").ColumnWidth = 12.166666666666666

# Headers + new "unit" values.
$wsProd.Range("D1").Value = "unit"
$wsProd.Range("D2").Value = "kg"
$wsProd.Range("D3").Value = "pc"

# New "discountType" header + values.
$wsProd.Range("I1").Value = "discountType"
$wsProd.Range("I3").Value = "Flat"
$wsProd.Range("I2").Value = "Percent"

# catName now points at the demo category for both sample rows.
$wsProd.Range("B2").Value = "Demo category 1"
$wsProd.Range("B3").Value = "Demo category 1"

$null = $wsProd.Range("B5").Select()
